# Update cryptocurrency price/volume data as of the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.210.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7017"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.82%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07809"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.71%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3108"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

# Row 10
$ws.Range("E10").Value = "  -4.19%  "

# Row 11
$ws.Range("E11").Value = "  -3.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "92.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.117"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6884"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.523"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008427"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "

# Row 18
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.097.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "

# Row 21
$ws.Range("E21").Value = "  -3.18%  "

# Row 22
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.587"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1533"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.879"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.92%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.567"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.02%  "

# Row 31
$ws.Range("E31").Value = "  -1.45%  "

# Row 32
$ws.Range("E32").Value = "  -1.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05215"
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7582"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.871"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.37%  "

# Row 38
$ws.Range("E38").Value = "  -1.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.225.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.721"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8981"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.748"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.04%  "

# Row 44
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.004.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.64%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.75%  "

# Row 48
$ws.Range("E48").Value = "  -0.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("E50").Value = "  -2.48%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.026"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
